$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.022.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.23%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.869.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.27%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.19%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'311.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.41%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.5159"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.35%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.3848"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.16%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.08276"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.60%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "'  -0.37%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'41.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.54%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'6.209"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.28%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.880.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.27%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "'20.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.46%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'7.307"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.34%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.18%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.00001099"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.06%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'90.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.21%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.06637"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.18%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'17.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.96%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  +0.09%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'6.034"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.02%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'28.056.33"
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'11.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.16%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "'  -0.61%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'2.081.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.49%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'2.501"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.47%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'157.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.37%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'20.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.78%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  -1.00%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'0.1067"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.86%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  -2.78%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'5.811"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.51%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'3.595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.18%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'9.490"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.65%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  -1.27%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'0.06507"
$ws.Range("D37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.2204"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.04%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.6591"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.18%  "
$ws.Range("E39").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'5.018"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.41%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  -2.70%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'11.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.31%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.6146"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.28%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'13.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.17%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "'  +0.13%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'3.665"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.03%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'2.028"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.82%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'1.219"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.21%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'120.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.45%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'78.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.54%  "
$ws.Range("E51").Style = "Normal"
